# Add a new "Leave Card" entry (VL(1-0-0) and a second VL(2-0-0) occurrence)
# to the leave table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# --- Row 50 (period 7/1/2023, SL(1-0-0)) : fill in EARNED days ---
$ws.Range("C50").Value = 1.25

# --- Row 51 (period 8/1/2023) : record a VL(1-0-0) leave credit ---
$ws.Range("B51").Value = "VL(1-0-0)"
$ws.Range("C51").Value = 1.25
$ws.Range("D51").Value = 1
# Copy the date format from the row above before writing the date value
$ws.Range("K50").Copy()
$ws.Range("K51").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K51").Value = 45156

# --- Row 52 (period 9/1/2023) : first VL(2-0-0) occurrence in September ---
$ws.Range("B52").Value = "VL(2-0-0)"
$ws.Range("D52").Value = 2
$ws.Range("K52").Value = "9/7,14/2023"

# --- Insert a new table row (53) for the second VL(2-0-0) occurrence in September ---
$ws.Rows.Item(53).Insert()
$ws.Range("A54:K54").Copy()
$ws.Range("A53:K53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Resize (extend) the table to include the newly inserted row before touching
# any calculated-column formulas so structured references resolve correctly.
$tbl.Resize($ws.Range("A8:K132"))

# Refresh the calculated "EARNED " column formula for the rows affected by the insert.
$ws.Range("G53").Formula = $ws.Range("G54").Formula()
$ws.Range("G132").Formula = $ws.Range("G132").Formula()

# Fill in the new row's data (A53 stays blank - same September period as row 52).
$ws.Range("B53").Value = "VL(2-0-0)"
$ws.Range("D53").Value = 2
$ws.Range("K53").Value = "9/21,28/2023"

# Recalculate so BALANCE formulas (E9 / I9, etc.) refresh with the new totals.
$excel.CalculateFull()

# Restore the last active-cell selection recorded in the saved file.
[void]$ws.Activate()
[void]$ws.Range("E53").Select()
